# Weekly update: insert a new data row for the latest week at row 81,
# pushing the existing historical rows (81-102) down to (82-103).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81 (shifts rows 81:102 down to 82:103,
# copies formatting from the row above, same as native Excel behaviour).
$row = $ws.Rows.Item(81)
$row.Insert()

# Populate the newly inserted row 81 with this week's record.
$ws.Cells.Item(81, 1).Value = 3
$ws.Cells.Item(81, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(81, 3).Value = "Coquimbo"
$ws.Cells.Item(81, 4).Value = 44551
$ws.Cells.Item(81, 5).Value = 5
$ws.Cells.Item(81, 6).Value = 100112052
$ws.Cells.Item(81, 7).Value = "Albahaca"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 105
$ws.Cells.Item(81, 11).Value = 4500
$ws.Cells.Item(81, 12).Value = 5000
$ws.Cells.Item(81, 13).Value = 4738
$ws.Cells.Item(81, 14).Value = "$/docena de matas"
$ws.Cells.Item(81, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(81, 16).Value = 790
$ws.Cells.Item(81, 17).Value = 6
$ws.Cells.Item(81, 18).Value = "Hortaliza"
